# Insert a new row at 193, shifting existing data (rows 193-280) down to
# rows 194-281. This matches the diff: every existing row from 193 onward
# moves down by one row, the former last row (280) becomes row 281, and a
# brand new row of data is written into the now-empty row 193. The sheet
# dimension grows from A1:R280 to A1:R281 automatically as part of the
# insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new record's data.
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 44636
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = 100114014
$ws.Range("G193").Value = "Betarraga"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 400
$ws.Range("K193").Value = 600
$ws.Range("L193").Value = 650
$ws.Range("M193").Value = 625
$ws.Range("N193").Value = '$/paquete 5 unidades'
$ws.Range("O193").Value = "Región Metropolitana"
$ws.Range("P193").Value = 125
$ws.Range("Q193").Value = 5
$ws.Range("R193").Value = "Hortaliza"
